# regen save_data to use K (actual strikeouts) instead of Strike# (strike-pitch
# count), regen std/mean, calc and write s_vals.
#
# Column G (header "K") previously held a "Strike#" style value (count of
# strikes thrown during the outing). It is being replaced here with the
# actual strikeout total (K) pulled from the regenerated source data for
# each start of the season.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> new K value (column G), keyed by the row's position in the
# worksheet (row 2 = most recent start, row 46 = earliest start in the log).
$kValues = @{
    2  = 1
    3  = 3
    4  = 0
    5  = 2
    6  = 0
    7  = 3
    8  = 1
    9  = 2
    10 = 2
    11 = 2
    12 = 0
    13 = 2
    14 = 1
    15 = 0
    17 = 1
    18 = 2
    19 = 1
    20 = 0
    21 = 2
    22 = 1
    23 = 1
    24 = 3
    25 = 2
    26 = 2
    27 = 3
    28 = 0
    29 = 4
    30 = 0
    31 = 2
    32 = 2
    33 = 1
    34 = 1
    35 = 0
    36 = 2
    37 = 2
    38 = 1
    39 = 1
    40 = 7
    41 = 1
    42 = 1
    44 = 1
    45 = 2
    46 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
